# Catalog SQL bugs: replace the CERT-code rows (A10:B20) with the new
# CWE-code catalog (A10:B41), extending the table from 18 data rows to
# 32 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 and 11 get new labels/values (same shared-string "slot" reused
# for CWE-126 / CWE-823 replacing the old MEM35-C / ARR30-C values).
$ws.Range("A10").Value = "CWE-126"
$ws.Range("B10").Value = 62

$ws.Range("A11").Value = "CWE-823"
$ws.Range("B11").Value = 42

# Row 12's value (1) is unchanged, but its label text changes from the
# old CERT code to the new CWE code.
$ws.Range("A12").Value = "CWE-196"
$ws.Range("B12").Value = 1

# Rows 13 and 14 swap label order relative to the old sheet.
$ws.Range("A13").Value = "CWE-754"
$ws.Range("B13").Value = 7

$ws.Range("A14").Value = "CWE-824"
$ws.Range("B14").Value = 26

$ws.Range("A15").Value = "CWE-457"
$ws.Range("B15").Value = 9

$ws.Range("A16").Value = "CWE-416"
$ws.Range("B16").Value = 13

$ws.Range("A17").Value = "CWE-119"
$ws.Range("B17").Value = 54

$ws.Range("A18").Value = "CWE-835"
$ws.Range("B18").Value = 30

$ws.Range("A19").Value = "CWE-822"
$ws.Range("B19").Value = 6

$ws.Range("A20").Value = "CWE-125"
$ws.Range("B20").Value = 10

# New rows 21-41 extending the CWE catalog.
$ws.Range("A21").Value = "CWE-469"
$ws.Range("B21").Value = 14

$ws.Range("A22").Value = "CWE-415"
$ws.Range("B22").Value = 11

$ws.Range("A23").Value = "CWE-476"
$ws.Range("B23").Value = 42

$ws.Range("A24").Value = "CWE-120"
$ws.Range("B24").Value = 32

$ws.Range("A25").Value = "CWE-834"
$ws.Range("B25").Value = 7

$ws.Range("A26").Value = "CWE-787"
$ws.Range("B26").Value = 1

$ws.Range("A27").Value = "CWE-191"
$ws.Range("B27").Value = 1

$ws.Range("A28").Value = "CWE-126"
$ws.Range("B28").Value = 2

$ws.Range("A29").Value = "CWE-788"
$ws.Range("B29").Value = 1

$ws.Range("A30").Value = "CWE-628"
$ws.Range("B30").Value = 1

$ws.Range("A31").Value = "CWE-131"
$ws.Range("B31").Value = 33

$ws.Range("A32").Value = "CWE-170"
$ws.Range("B32").Value = 2

$ws.Range("A33").Value = "CWE-129"
$ws.Range("B33").Value = 10

$ws.Range("A34").Value = "CWE-460"
$ws.Range("B34").Value = 5

$ws.Range("A35").Value = "CWE-825"
$ws.Range("B35").Value = 4

$ws.Range("A36").Value = "CWE-121"
$ws.Range("B36").Value = 6

$ws.Range("A37").Value = "CWE-190"
$ws.Range("B37").Value = 4

$ws.Range("A38").Value = "CWE-789"
$ws.Range("B38").Value = 4

$ws.Range("A39").Value = "CWE-248"
$ws.Range("B39").Value = 4

$ws.Range("A40").Value = "CWE-127"
$ws.Range("B40").Value = 4

$ws.Range("A41").Value = "CWE-124"
$ws.Range("B41").Value = 4

# Match the saved selection state from the diff (active cell A10,
# selection spans the full row A10:XFD10).
$ws.Range("A10:XFD10").Select()
